# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.623.46"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "2.594.46"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Formula = "'1.08"
$ws.Range("E4").Value = "  +8.08%  "
$ws.Range("D5").Formula = "'568.97"
$ws.Range("E5").Value = "  +4.96%  "
$ws.Range("D6").Formula = "'143.18"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("D9").Value = "2.601.51"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Formula = "'6.63"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("E12").Value = "  +9.92%  "
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "3.059.50"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "59.669.40"
$ws.Range("D16").Formula = "'21.89"
$ws.Range("E16").Value = "  +6.73%  "
$ws.Range("D17").Formula = "'0.0000136"
$ws.Range("E17").Value = "  +3.55%  "
$ws.Range("D18").Value = "2.596.67"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Formula = "'337.31"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Formula = "'10.22"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Formula = "'0.445"
$ws.Range("E25").Value = "  +6.38%  "
$ws.Range("D26").Formula = "'0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").Formula = "'7.27"
$ws.Range("E28").Value = "  +3.74%  "
$ws.Range("E29").Value = "  +7.56%  "
$ws.Range("D30").Formula = "'0.999"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("D32").Formula = "'161.06"
$ws.Range("E32").Value = "  +4.42%  "
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").Formula = "'18.98"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Formula = "'4.04"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("D36").Formula = "'0.890"
$ws.Range("E36").Value = "  +9.74%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Formula = "'0.876"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Formula = "'1.13"
$ws.Range("E38").Value = "  +4.34%  "
$ws.Range("D39").Formula = "'37.21"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D41").Formula = "'293.12"
$ws.Range("E41").Value = "  +5.53%  "
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("D43").Formula = "'0.998"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Formula = "'19.01"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("D49").Formula = "'125.21"
$ws.Range("E49").Value = "  +15.19%  "
$ws.Range("D50").Formula = "'0.0232"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").Value = "1.934.15"
$ws.Range("E51").Value = "  +1.70%  "
